$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "runs"/"balls"/"fours" (columns C/D/E) for the two Nathan Coulter-Nile
# innings rows were swapped - row 2 should now show row 3's figures and
# vice versa. The source values are stored as text (e.g. "24"), so swap
# them via Copy/PasteSpecial (through a scratch cell) rather than
# re-assigning .Value directly, which would let Excel's COM layer
# auto-convert the numeric-looking text into real numbers.
$temp = $ws.Range("Z1")

foreach ($col in @("C", "D", "E")) {
    $cellRow2 = $ws.Range($col + "2")
    $cellRow3 = $ws.Range($col + "3")

    $cellRow2.Copy()
    $temp.PasteSpecial()

    $cellRow3.Copy()
    $cellRow2.PasteSpecial()

    $temp.Copy()
    $cellRow3.PasteSpecial()

    $temp.Clear()
}
